$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.493.06"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.911.83"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "325.17"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").Value = "0.4061"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "1.010"
$ws.Range("D11").Value = "23.45"
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("D12").Value = "1.920.05"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "5.981"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "7.114"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "90.11"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "0.06770"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "17.66"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "29.511.68"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "5.614"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").Value = "2.181"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "2.138.65"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "155.25"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "6.369"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "2.100"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "119.84"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").Value = "0.09524"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "5.497"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "3.567"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "1.385"
$ws.Range("E35").Value = "  -3.32%  "
$ws.Range("D36").Value = "0.02265"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "0.06097"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "0.5929"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "7.945"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Value = "0.1851"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").Value = "2.412"
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").Value = "0.07643"
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("D46").Value = "12.51"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").Value = "0.5570"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "116.19"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").Value = "72.36"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").Value = "  +1.80%  "
